# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new trailing columns (date, legislator_name, legislator_id) to the
# "股票" (stock) sheet, populating the header row and the three data rows with
# the legislator's filing date, name, and id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- Header row (row 1): copy the existing bold/bordered header format from
# G1 so the new header cells (H1:J1) render identically to the rest of the
# header row, then fill in the new column names.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# ---- Data rows (2-4): the "date" column holds a literal yyyy-mm-dd string,
# not a real date -- format the cells as Text first so Excel doesn't
# auto-convert "2011-11-18" into a date serial number.
$ws.Range("H2:H4").NumberFormat = "@"
$ws.Range("H2").Value = "2011-11-18"
$ws.Range("H3").Value = "2011-11-18"
$ws.Range("H4").Value = "2011-11-18"

$ws.Range("I2").Value = "張嘉郡"
$ws.Range("I3").Value = "張嘉郡"
$ws.Range("I4").Value = "張嘉郡"

$ws.Range("J2").Value = 1719
$ws.Range("J3").Value = 1719
$ws.Range("J4").Value = 1719

# Clear the temporary Text number format again so the data cells fall back to
# the sheet's normal (unstyled) look, matching the rest of the table.
$ws.Range("H2:J4").Style = "Normal"
